# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# handed off for both zh-cn and de-de locales: new handoff xlf files were
# generated, the status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and an error detail notes the handback file is stale.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/db61008db46db3c80c09c7b7a35f6f28c62b2ba0/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3beb121403298d9d8446def1217a7a3cb8f9557a/e2e/b.md."

# ---- Overview sheet: summary row for b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-26 08:41:47"

# ---- zh-cn sheet: detail row for b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-26 08:41:43"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet: detail row for b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-26 08:41:47"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
